{"js": "// Apply the benchmark-stats update to the single-column results table.\n// Each table row holds one cell; we rewrite the cell text for the rows\n// whose values changed, by row index (0-based), and collapse the three\n// multi-run (tab-separated) summary rows near the end of the table into\n// plain single-value cells.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of table row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"680\",\n  4: \"0.00002\",\n  6: \"0.00012\",\n  8: \"0.00021\",\n  9: \"0.00023\",\n  10: \"0.00026\",\n  11: \"0.09338\",\n  43: \"99.88\",\n  44: \"0.09\",\n  45: \"76\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(parseInt(rowIndex, 10), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stats update to the single-column results table.\n# Each table row holds one cell; we rewrite the cell text for the rows\n# whose values changed (1-based row indices, as used by Word's COM\n# Table.Cell(row, col)), collapsing the three multi-run (tab-separated)\n# summary rows near the end of the table into plain single-value cells.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"680\"\n    5  = \"0.00002\"\n    7  = \"0.00012\"\n    9  = \"0.00021\"\n    10 = \"0.00023\"\n    11 = \"0.00026\"\n    12 = \"0.09338\"\n    44 = \"99.88\"\n    45 = \"0.09\"\n    46 = \"76\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $t.Cell($rowIndex, 1).Range.Text = $updates[$rowIndex]\n}\n"}
